# Apply "Commiting changes of email reports" edit to the Queries sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Queries")

# xlVAlignTop
$xlTop = -4160

# --- 1. Insert a new column before A for "DB Type" (shifts B..G -> C..H,
#         preserving their widths/styles/values exactly). -----------------
$ws.Columns("A:A").Insert()

# --- 2. Populate the new column A header + row2 value, with the same
#         "vertical top" cell style (style index 4) used by the rest of
#         row 2. -----------------------------------------------------------
$ws.Range("A1").Value2 = "DB Type"
$ws.Range("A1").VerticalAlignment = $xlTop

$ws.Range("A2").Value2 = "MS SQL"
$ws.Range("A2").VerticalAlignment = $xlTop

# --- 3. Column widths: new column A + the newly-explicit column B
#         (formerly default-width "Report Channel" column). ---------------
$ws.Columns("A:A").ColumnWidth = 8.7265625
$ws.Columns("B:B").ColumnWidth = 13.81640625

# --- 4. Add a new row 3 (MYSQL variant) by copying row 2 (carries over all
#         per-cell styles/formats) then overwriting the cells that differ. -
$ws.Rows("2:2").Copy()
$ws.Rows("3:3").Insert()

$ws.Range("A3").Value2 = "MYSQL"
# Leading apostrophe keeps the "quotePrefix" (style 3) formatting these
# date-look-alike strings already had (matches E2/F2), while the engine
# strips the apostrophe itself from the stored text.
$ws.Range("E3").Value2 = "'01-06-2021 00:00:00"
$ws.Range("F3").Value2 = "'07-06-2021 00:00:00"
$ws.Range("G3").Value2 = "Select * from AGT_Agent_TimeTrack;"

$ws.Rows("3:3").RowHeight = 406

# --- 5. Update the active selection to match the authored state. ----------
$ws.Range("G3").Select()
